$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. The "_GoBack" bookmark currently sits right after the trailing
#    space in the "Reform mat the code " paragraph. In the edited
#    document it has moved further down (into the new "Week 5"
#    content), so remove it from its old spot first.
# ------------------------------------------------------------------
try {
    $oldBookmark = $d.Bookmarks.Item("_GoBack")
    $oldBookmark.Delete()
} catch {
    # no pre-existing bookmark -- nothing to remove
}

# ------------------------------------------------------------------
# 2. Append the "Week 3 summary forgot to commit" content: new
#    "Find a better way..." / "Make a grammar tree" / "Make it more
#    faster to run" bullets, a new "Week 5" heading, and the two
#    closing paragraphs -- all inserted right after the existing
#    "Week 4" heading paragraph.
# ------------------------------------------------------------------
$weekFour = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.TrimEnd([char]13, [char]7) -eq "Week 4") {
        $weekFour = $candidate
    }
}

$insertionPoint = $weekFour.Range
$insertionPoint.Collapse(0)
$insertionPoint.InsertParagraphAfter()
$newFirstParaIndex = $weekFour.Index + 1
$target = $d.Paragraphs.Item($newFirstParaIndex).Range

$xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:pPr><w:rPr><w:highlight w:val="green"/></w:rPr></w:pPr><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>Find a better way to use parser</w:t></w:r></w:p>
<w:p><w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t>Make a grammar tree</w:t></w:r></w:p>
<w:p><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">Make it </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>more faster</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> to run</w:t></w:r></w:p>
<w:p><w:pPr><w:rPr><w:i/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:u w:val="single"/></w:rPr><w:t>Week 5</w:t></w:r></w:p>
<w:p><w:r><w:lastRenderedPageBreak/><w:t>Using the basic grammar tree,</w:t></w:r></w:p>
<w:p><w:r><w:t>Make at least 10+</w:t></w:r></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$null = $target.InsertXML($xml)

# ------------------------------------------------------------------
# 3. Re-add the "_GoBack" bookmark inside the new
#    "Using the basic grammar tree," paragraph, splitting it right
#    after "...grammar t" / before "ree,".
# ------------------------------------------------------------------
$grammarParaIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text.StartsWith("Using the basic grammar tree,")) {
        $grammarParaIndex = $i
    }
}

$grammarPara = $d.Paragraphs.Item($grammarParaIndex)
$splitOffset = $grammarPara.Range.Start + "Using the basic grammar t".Length
$bookmarkRange = $d.Range($splitOffset, $splitOffset)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
